$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "mse" header in I1, matching the formatting of the other header
# cells (bold, bordered, centered) by copying H1's format onto I1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "mse"

# New mse data values for rows 2 and 3
$ws.Range("I2").Value = 0.8011432342852155
$ws.Range("I3").Value = 0.7594823141892751

# Updated values for row 3 (regression_binary_pred / F_1_score)
$ws.Range("G3").Value = 0.5937239583333334
$ws.Range("H3").Value = 0.7365552616033121
